$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "processor_full_name" column (column A) entirely, shifting
# all subsequent columns one to the left.
$ws.Range("A1").EntireColumn.Delete()

$ws.Range("A1").Select()
